# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Estados Unidos (row 4) - updated case counts
$ws.Range("B4").Value = 1131856
$ws.Range("C4").Value = 826
$ws.Range("D4").Value = 161666
$ws.Range("E4").Value = 904408
$ws.Range("F4").Value = 16481
$ws.Range("G4").Value = 29
$ws.Range("H4").Value = 65782

# Reorder Catar / Japon (rows 33-34) and refresh their figures
$ws.Range("A33").Value = "Catar"
$ws.Range("B33").Value = 14872
$ws.Range("C33").Value = 776
$ws.Range("D33").Value = 1534
$ws.Range("E33").Value = 13326
$ws.Range("F33").Value = 72
$ws.Range("G33").Value = 0
$ws.Range("H33").Value = 12

$ws.Range("A34").Value = "Japon"
$ws.Range("B34").Value = 14305
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 2975
$ws.Range("E34").Value = 10875
$ws.Range("F34").Value = 328
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 455

# Eslovenia (row 84) - updated case counts
$ws.Range("B84").Value = 1439
$ws.Range("C84").Value = 5
$ws.Range("D84").Value = 239
$ws.Range("E84").Value = 1106
$ws.Range("F84").Value = 21
$ws.Range("G84").Value = 2
$ws.Range("H84").Value = 94

# Republica de Yibuti (row 89) - updated case counts
$ws.Range("B89").Value = 1112
$ws.Range("C89").Value = 15
$ws.Range("D89").Value = 686
$ws.Range("E89").Value = 424

# Malta (row 116) - updated case counts
$ws.Range("B116").Value = 468
$ws.Range("C116").Value = 1
$ws.Range("D116").Value = 379
$ws.Range("E116").Value = 85

# Reorder Zambia / Trinidad y Tobago / Bermudas (rows 149-151) and refresh Zambia's figures
$ws.Range("A149").Value = "Zambia"
$ws.Range("B149").Value = 119
$ws.Range("C149").Value = 10
$ws.Range("D149").Value = 75
$ws.Range("E149").Value = 41
$ws.Range("F149").Value = 1
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 3

$ws.Range("A150").Value = "Trinidad yTobago"
$ws.Range("B150").Value = 116
$ws.Range("C150").Value = 0
$ws.Range("D150").Value = 83
$ws.Range("E150").Value = 25
$ws.Range("F150").Value = 0
$ws.Range("G150").Value = 0
$ws.Range("H150").Value = 8

$ws.Range("A151").Value = "Bermudas"
$ws.Range("B151").Value = 114
$ws.Range("C151").Value = 0
$ws.Range("D151").Value = 48
$ws.Range("E151").Value = 60
$ws.Range("F151").Value = 4
$ws.Range("G151").Value = 0
$ws.Range("H151").Value = 6

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 2 de Mayo de 2020 a las 12:48"
